{"js": "// Replace the worksheet date and every \"A\u00f7B=C, D\" answer cell with the\n// updated values from the commit. Each old string is unique in the\n// document, so a plain-text, case-sensitive search + full-text replace\n// on each hit is sufficient and avoids relying on table/row ordering.\nconst replacements = [\n  [\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"],\n  [\"315\u00f75=63, 0\", \"345\u00f73=115, 0\"],\n  [\"342\u00f77=48, 6\", \"160\u00f74=40, 0\"],\n  [\"860\u00f78=107, 4\", \"152\u00f73=50, 2\"],\n  [\"622\u00f75=124, 2\", \"654\u00f76=109, 0\"],\n  [\"865\u00f72=432, 1\", \"257\u00f72=128, 1\"],\n  [\"902\u00f79=100, 2\", \"215\u00f75=43, 0\"],\n  [\"302\u00f78=37, 6\", \"184\u00f79=20, 4\"],\n  [\"614\u00f79=68, 2\", \"828\u00f78=103, 4\"],\n  [\"131\u00f75=26, 1\", \"601\u00f77=85, 6\"],\n  [\"739\u00f72=369, 1\", \"900\u00f75=180, 0\"],\n  [\"869\u00f77=124, 1\", \"835\u00f78=104, 3\"],\n  [\"723\u00f74=180, 3\", \"713\u00f79=79, 2\"],\n  [\"557\u00f76=92, 5\", \"551\u00f79=61, 2\"],\n  [\"325\u00f79=36, 1\", \"422\u00f79=46, 8\"],\n  [\"536\u00f75=107, 1\", \"778\u00f72=389, 0\"],\n  [\"811\u00f73=270, 1\", \"398\u00f78=49, 6\"],\n  [\"326\u00f75=65, 1\", \"917\u00f77=131, 0\"],\n  [\"376\u00f77=53, 5\", \"527\u00f72=263, 1\"],\n  [\"918\u00f75=183, 3\", \"471\u00f75=94, 1\"],\n  [\"107\u00f76=17, 5\", \"448\u00f74=112, 0\"],\n  [\"698\u00f73=232, 2\", \"744\u00f76=124, 0\"],\n  [\"231\u00f79=25, 6\", \"414\u00f73=138, 0\"],\n  [\"443\u00f74=110, 3\", \"405\u00f74=101, 1\"],\n  [\"805\u00f78=100, 5\", \"940\u00f74=235, 0\"],\n  [\"223\u00f74=55, 3\", \"116\u00f74=29, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the worksheet date and every \"A\u00f7B=C, D\" answer cell with the\n# updated values from the commit. Each old string is unique in the\n# document, so Find/Replace (wdReplaceOne) scoped to the whole story via\n# $d.Content is sufficient and does not depend on table/row ordering.\n$pairs = @(\n  @(\"2024-01-20 Saturday\", \"2024-01-21 Sunday\"),\n  @(\"315\u00f75=63, 0\", \"345\u00f73=115, 0\"),\n  @(\"342\u00f77=48, 6\", \"160\u00f74=40, 0\"),\n  @(\"860\u00f78=107, 4\", \"152\u00f73=50, 2\"),\n  @(\"622\u00f75=124, 2\", \"654\u00f76=109, 0\"),\n  @(\"865\u00f72=432, 1\", \"257\u00f72=128, 1\"),\n  @(\"902\u00f79=100, 2\", \"215\u00f75=43, 0\"),\n  @(\"302\u00f78=37, 6\", \"184\u00f79=20, 4\"),\n  @(\"614\u00f79=68, 2\", \"828\u00f78=103, 4\"),\n  @(\"131\u00f75=26, 1\", \"601\u00f77=85, 6\"),\n  @(\"739\u00f72=369, 1\", \"900\u00f75=180, 0\"),\n  @(\"869\u00f77=124, 1\", \"835\u00f78=104, 3\"),\n  @(\"723\u00f74=180, 3\", \"713\u00f79=79, 2\"),\n  @(\"557\u00f76=92, 5\", \"551\u00f79=61, 2\"),\n  @(\"325\u00f79=36, 1\", \"422\u00f79=46, 8\"),\n  @(\"536\u00f75=107, 1\", \"778\u00f72=389, 0\"),\n  @(\"811\u00f73=270, 1\", \"398\u00f78=49, 6\"),\n  @(\"326\u00f75=65, 1\", \"917\u00f77=131, 0\"),\n  @(\"376\u00f77=53, 5\", \"527\u00f72=263, 1\"),\n  @(\"918\u00f75=183, 3\", \"471\u00f75=94, 1\"),\n  @(\"107\u00f76=17, 5\", \"448\u00f74=112, 0\"),\n  @(\"698\u00f73=232, 2\", \"744\u00f76=124, 0\"),\n  @(\"231\u00f79=25, 6\", \"414\u00f73=138, 0\"),\n  @(\"443\u00f74=110, 3\", \"405\u00f74=101, 1\"),\n  @(\"805\u00f78=100, 5\", \"940\u00f74=235, 0\"),\n  @(\"223\u00f74=55, 3\", \"116\u00f74=29, 0\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
